# Insert a new "Control Signals" section (3 blank paragraphs + heading)
# right after the existing "FSM Diagram" paragraph, mirroring the formatting
# of that paragraph (ListParagraph style, centered, Cambria, bold, sz 180).

$d = $word.ActiveDocument

# Locate the "FSM Diagram" paragraph via Find.
$range = $d.Content
$found = $range.Find.Execute("FSM Diagram", $false, $false, $false, $false, $false,
                              $true, 1, $false, "", 0)

# Collapse to the end of the found text, then move to cover the whole paragraph
# so we can grab its Range for InsertParagraphAfter.
$range.Collapse(0)  # wdCollapseEnd

# Insert four new paragraphs after the "FSM Diagram" paragraph.
$p1 = $range.InsertParagraphAfter()
$range.Collapse(0)

$p2 = $range.InsertParagraphAfter()
$range.Collapse(0)

$p3 = $range.InsertParagraphAfter()
$range.Collapse(0)

$p4 = $range.InsertParagraphAfter()
$range.Collapse(0)

# The four newly inserted paragraphs are now the next four paragraphs after
# the "FSM Diagram" one. Grab them via the Paragraphs collection.
$fsmPara = $range.Paragraphs(1)

# Re-find the FSM Diagram paragraph's index so we can address the following ones.
$allParas = $d.Paragraphs
$fsmIndex = -1
for ($i = 1; $i -le $allParas.Count; $i++) {
    if ($allParas.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "FSM Diagram") {
        $fsmIndex = $i
    }
}

for ($offset = 1; $offset -le 4; $offset++) {
    $para = $allParas.Item($fsmIndex + $offset)
    $pr = $para.Range
    $pr.ParagraphFormat.Style = "List Paragraph"
    $pr.ParagraphFormat.Alignment = 1  # wdAlignParagraphCenter
    $pr.Font.Name = "Cambria"
    $pr.Font.Bold = $true
    $pr.Font.Size = 90
}

# Put the heading text on the 4th of the newly added paragraphs (the last one).
$lastPara = $allParas.Item($fsmIndex + 4)
$lastRange = $lastPara.Range
$lastRange.Collapse(0)
$lastRange.MoveEnd(1, -1) | Out-Null
$lastRange.Text = "Control Signals"
$lastRange.Font.Name = "Cambria"
$lastRange.Font.Bold = $true
$lastRange.Font.Size = 90
